$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to re-pulled data per commit message:
# "repull data, push all data, mean calculation"
$ws.Range("F3").Value = -11
$ws.Range("F4").Value = -5
$ws.Range("F5").Value = 9
$ws.Range("F8").Value = 13
$ws.Range("F10").Value = -3
